$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("T2").Value = 0.1450048780487805
$ws.Range("V2").Value = 0.0002448603057459146
$ws.Range("Z2").Value = -0.1844856619543747
$ws.Range("AB2").Value = -753.4322943540341
$ws.Range("AC2").Value = "umolO2/min/m2"
$ws.Range("AD2").Value = -753.4322943540341

# Row 3
$ws.Range("T3").Value = 0.1492487804878049
$ws.Range("V3").Value = 0.0001488973818309612
$ws.Range("Z3").Value = -0.1899027326049314
$ws.Range("AB3").Value = -1275.393363333429
$ws.Range("AC3").Value = "umolO2/min/m2"
$ws.Range("AD3").Value = -1275.393363333429

# Row 4
$ws.Range("T4").Value = 0.1469268292682927
$ws.Range("V4").Value = 0.0002222807942365138
$ws.Range("Z4").Value = -0.1808977077723057
$ws.Range("AB4").Value = -813.8251817645781
$ws.Range("AC4").Value = "umolO2/min/m2"
$ws.Range("AD4").Value = -813.8251817645781

# Row 5
$ws.Range("T5").Value = 0.1418926829268293
$ws.Range("V5").Value = 0.0002529432437181515
$ws.Range("Z5").Value = -0.2007310390331157
$ws.Range("AB5").Value = -793.5813429228629
$ws.Range("AC5").Value = "umolO2/min/m2"
$ws.Range("AD5").Value = -793.5813429228629

# Row 6
$ws.Range("T6").Value = 0.1446439024390244
$ws.Range("V6").Value = 0.0001851607801792304
$ws.Range("Z6").Value = -0.1444995990480928
$ws.Range("AB6").Value = -780.400681549415
$ws.Range("AC6").Value = "umolO2/min/m2"
$ws.Range("AD6").Value = -780.400681549415

# Row 7
$ws.Range("T7").Value = 0.1429268292682927
$ws.Range("V7").Value = 0.0003232296608680373
$ws.Range("Z7").Value = -0.1889044042065035
$ws.Range("AB7").Value = -584.427814264317
$ws.Range("AC7").Value = "umolO2/min/m2"
$ws.Range("AD7").Value = -584.427814264317

# Row 8
$ws.Range("T8").Value = 0.1544
$ws.Range("V8").Value = 0
$ws.Range("Z8").Value = 0.001263385345290537
$ws.Range("AB8").Value = "Inf"
$ws.Range("AC8").Value = "umolO2/min/m2"
$ws.Range("AD8").Value = "Inf"

# Row 9
$ws.Range("T9").Value = 0.1450048780487805
$ws.Range("V9").Value = 0.0002448603057459146
$ws.Range("Z9").Value = 0.05647306839756049
$ws.Range("AB9").Value = 230.6338229282503
$ws.Range("AC9").Value = "umolO2/min/m2"
$ws.Range("AD9").Value = 230.6338229282503

# Row 10
$ws.Range("T10").Value = 0.1492487804878049
$ws.Range("V10").Value = 0.0001488973818309612
$ws.Range("Z10").Value = 0.1174595621853158
$ws.Range("AB10").Value = 788.8625088026342
$ws.Range("AC10").Value = "umolO2/min/m2"
$ws.Range("AD10").Value = 788.8625088026342

# Row 11
$ws.Range("T11").Value = 0.1469268292682927
$ws.Range("V11").Value = 0.0002222807942365138
$ws.Range("Z11").Value = -0.02662312056838054
$ws.Range("AB11").Value = -119.7724736400424
$ws.Range("AC11").Value = "umolO2/min/m2"
$ws.Range("AD11").Value = -119.7724736400424

# Row 12
$ws.Range("T12").Value = 0.1418926829268293
$ws.Range("V12").Value = 0.0002529432437181515
$ws.Range("Z12").Value = 0.07428576836083688
$ws.Range("AB12").Value = 293.6855211820234
$ws.Range("AC12").Value = "umolO2/min/m2"
$ws.Range("AD12").Value = 293.6855211820234

# Row 13
$ws.Range("T13").Value = 0.1446439024390244
$ws.Range("V13").Value = 0.0001851607801792304
$ws.Range("Z13").Value = 0.1087342486918217
$ws.Range("AB13").Value = 587.2423338601731
$ws.Range("AC13").Value = "umolO2/min/m2"
$ws.Range("AD13").Value = 587.2423338601731

# Row 14
$ws.Range("T14").Value = 0.1429268292682927
$ws.Range("V14").Value = 0.0003232296608680373
$ws.Range("Z14").Value = 0.07734450793189028
$ws.Range("AB14").Value = 239.2865423432386
$ws.Range("AC14").Value = "umolO2/min/m2"
$ws.Range("AD14").Value = 239.2865423432386

# Row 15
$ws.Range("T15").Value = 0.1544
$ws.Range("V15").Value = 0
$ws.Range("Z15").Value = -0.0002791541157665243
$ws.Range("AB15").Value = "-Inf"
$ws.Range("AC15").Value = "umolO2/min/m2"
$ws.Range("AD15").Value = "-Inf"
